$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for Fertilizer Recommendation
$ws.Range("U1").Value = "Fertilizer Recommendation"

# Update row 2 data values
$ws.Range("C2").Value = 656
$ws.Range("D2").Value = 656
$ws.Range("E2").Value = "asdasdasd"
$ws.Range("F2").Value = 52
$ws.Range("H2").Value = 65
$ws.Range("I2").Value = "23323as2d3asd"
$ws.Range("J2").Value = "6556456464"
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 232
$ws.Range("N2").Value = 200
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 21
$ws.Range("R2").Value = 23
$ws.Range("S2").Value = 0.331370118203286

# New Fertilizer Recommendation value for row 2
$ws.Range("U2").Value = "No specific fertilizer recommendation available for the given soil data."
